# Apply the "gh-pages output generated at 456a3b4" update to
# the 江西-漫展信息 workbook.
#
# Sheet 1 "展览"      -> rows with updated "想去人数" counts, the 九江·动漫畅想
#                        event being marked cancelled / not sellable.
# Sheet 2 "演出"      -> updated show time, attendee count and cover image.
# Sheet 3 "本地生活"  -> no changes (header row only).
# Sheet 4 "全部类型"  -> mirrors sheet 1 + sheet 2 changes (offset by the extra
#                        "演出" row that sheet 4 contains).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览"
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

$wsExpo.Range("F3").Value  = 5232
$wsExpo.Range("F5").Value  = 7531
$wsExpo.Range("F7").Value  = 74
$wsExpo.Range("F8").Value  = 102
$wsExpo.Range("F12").Value = 4353
$wsExpo.Range("F13").Value = 1774
$wsExpo.Range("F15").Value = 111
$wsExpo.Range("F16").Value = 2940
$wsExpo.Range("F20").Value = 519
$wsExpo.Range("F21").Value = 454
$wsExpo.Range("F23").Value = 321
$wsExpo.Range("F24").Value = 108
$wsExpo.Range("F25").Value = 1703
$wsExpo.Range("F26").Value = 1203
$wsExpo.Range("F28").Value = 1399

# Event cancelled - name suffixed, lowest-price column becomes "not sellable"
$wsExpo.Range("C33").Value = "九江·动漫畅想（取消）"
$wsExpo.Range("G33").Value = "不可售"

$wsExpo.Range("F34").Value = 9
$wsExpo.Range("F35").Value = 65
$wsExpo.Range("F37").Value = 71
$wsExpo.Range("F38").Value = 2970
$wsExpo.Range("F40").Value = 30
$wsExpo.Range("F41").Value = 101
$wsExpo.Range("F43").Value = 56

# ---------------------------------------------------------------------
# Sheet "演出"
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")

$wsShow.Range("E3").Value = "2024.08.17 14:00-08.17 15:30"
$wsShow.Range("F3").Value = 16
$wsShow.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202407/t07f8mmz1721894837940.jpeg"

# ---------------------------------------------------------------------
# Sheet "全部类型" (combined listing - same updates, rows shifted by +1
# starting row 23 because it also contains the "演出" entry as row 23/40)
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Range("F3").Value  = 5232
$wsAll.Range("F5").Value  = 7531
$wsAll.Range("F7").Value  = 74
$wsAll.Range("F8").Value  = 102
$wsAll.Range("F12").Value = 4353
$wsAll.Range("F13").Value = 1774
$wsAll.Range("F15").Value = 111
$wsAll.Range("F16").Value = 2940
$wsAll.Range("F20").Value = 519
$wsAll.Range("F21").Value = 454
$wsAll.Range("F24").Value = 321
$wsAll.Range("F25").Value = 108
$wsAll.Range("F26").Value = 1703
$wsAll.Range("F27").Value = 1203
$wsAll.Range("F29").Value = 1399

$wsAll.Range("C34").Value = "九江·动漫畅想（取消）"
$wsAll.Range("G34").Value = "不可售"

$wsAll.Range("F35").Value = 9
$wsAll.Range("F36").Value = 65
$wsAll.Range("F38").Value = 71
$wsAll.Range("F39").Value = 2970

$wsAll.Range("E40").Value = "2024.08.17 14:00-08.17 15:30"
$wsAll.Range("F40").Value = 16
$wsAll.Range("I40").Value = "//i0.hdslb.com/bfs/openplatform/202407/t07f8mmz1721894837940.jpeg"

$wsAll.Range("F42").Value = 30
$wsAll.Range("F43").Value = 101
$wsAll.Range("F45").Value = 56
